$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")
$c1 = $ws1.Range("C1")
try {
    $c1.Borders(9).LineStyle = 1
    Write-Host "a ok"
} catch { Write-Host "a failed: $_" }

try {
    $c1.Borders.Item(9).ColorIndex = 1
    Write-Host "b ok"
} catch { Write-Host "b failed: $_" }

try {
    [void]$ws1.Range("C1").Borders.Item(9)
    Write-Host "c ok"
} catch { Write-Host "c failed: $_" }
